$d = $word.ActiveDocument

# 1. Title: merge "Installation von " + "diNo" + " auf dem Server" into one run.
$d.Content.Find.Execute("Installation von diNo auf dem Server", $true, $false, $false, $false, $false, $true, 1, $false, "Installation von diNo auf dem Server", 2) | Out-Null

# 2. "diNo.bak" + " wiederherstellen" -> single run.
$d.Content.Find.Execute("diNo.bak wiederherstellen", $true, $false, $false, $false, $false, $true, 1, $false, "diNo.bak wiederherstellen", 2) | Out-Null

# 3. "Berechtigungen prüfen unter " + "diNo" + "/Sicherheit/Benutzer:" -> single run.
$d.Content.Find.Execute("Berechtigungen prüfen unter diNo/Sicherheit/Benutzer:", $true, $false, $false, $false, $false, $true, 1, $false, "Berechtigungen prüfen unter diNo/Sicherheit/Benutzer:", 2) | Out-Null

# 4. "SN\" + "sg" + "-" + "fb" + "-lehrer " -> single run.
$d.Content.Find.Execute("SN\sg-fb-lehrer ", $true, $false, $false, $false, $false, $true, 1, $false, "SN\sg-fb-lehrer ", 2) | Out-Null

# 5. ": Benutzer " + "sekretariat" + ",  Passwort: " -> single run.
$d.Content.Find.Execute(": Benutzer sekretariat,  Passwort: ", $true, $false, $false, $false, $false, $true, 1, $false, ": Benutzer sekretariat,  Passwort: ", 2) | Out-Null

# 6. "diNo87§" -> "diNo87" + "!" (two runs), and move the _GoBack bookmark here.
$d.Content.Find.Execute("diNo87§", $true, $false, $false, $false, $false, $true, 1, $false, "diNo87!", 2) | Out-Null

# 7. "Laufwerke freigeben:" + " für " + "dino" + "-aktuell und Notenbögen" -> two runs.
$d.Content.Find.Execute(" für dino-aktuell und Notenbögen", $true, $false, $false, $false, $false, $true, 1, $false, " für dino-aktuell und Notenbögen", 2) | Out-Null

# 8. "Copy-Job einrichten vom Verzeichnis bin auf " + "dino" + "-aktuell" -> single run.
$d.Content.Find.Execute("Copy-Job einrichten vom Verzeichnis bin auf dino-aktuell", $true, $false, $false, $false, $false, $true, 1, $false, "Copy-Job einrichten vom Verzeichnis bin auf dino-aktuell", 2) | Out-Null

# 9. "Copy-Job einrichten vom Verzeichnis " + "dino" + "-aktuell zum Verwaltungsserver (Hol-Job)" -> single run.
$d.Content.Find.Execute("Copy-Job einrichten vom Verzeichnis dino-aktuell zum Verwaltungsserver (Hol-Job)", $true, $false, $false, $false, $false, $true, 1, $false, "Copy-Job einrichten vom Verzeichnis dino-aktuell zum Verwaltungsserver (Hol-Job)", 2) | Out-Null

# 10. Move the _GoBack bookmark from the "Für Visual Studio..." paragraph to the
#     paragraph ending in "diNo87!" (right after the "!" run).
$bm = $d.Bookmarks.Item("_GoBack")
$bmRange = $bm.Range
$bm.Delete()

$target = $d.Content.Find
$target.Execute("diNo87!") | Out-Null
$rng = $d.Content.Duplicate
$rng.Start = $bmRange.Start
$rng.End = $bmRange.Start

$findRange = $d.Content.Duplicate
$findRange.Find.Execute("diNo87!") | Out-Null
$insertPoint = $findRange.Duplicate
$insertPoint.Collapse(0)
$d.Bookmarks.Add("_GoBack", $insertPoint) | Out-Null
